$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the "23/05/2023" bullet ("8:00- Ajuste y finalización de la primera
# carta.") paragraph robustly (wildcard match avoids accented-character
# comparison glitches) instead of trusting a hard-coded paragraph index.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$idx832 = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "8:00-*primera carta.*") {
        $idx832 = $i
        break
    }
}
if ($idx832 -eq -1) {
    throw "Could not locate the '8:00- Ajuste y finalización...' paragraph"
}

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# 1) Split "8:00- Ajuste y finalización de la primera carta." into the new
#    multi-run sentence about the homecards template / Jesús flipcards.
# ---------------------------------------------------------------------------
$body132 = @'
<w:p>
  <w:pPr>
    <w:ind w:left="360"/>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>8:00</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>-</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> 12:30</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> Ajuste y finalización de </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>la plantilla para el resto de las homecards, creación de las imágenes y los textos para los flipcards de la página Jesús.</w:t>
  </w:r>
</w:p>
'@

$range132 = $d.Paragraphs.Item($idx832).Range
$range132.InsertXML($pkgHeader + $body132 + $pkgFooter) | Out-Null

# ---------------------------------------------------------------------------
# 2) The paragraph right after it used to be an empty, bold placeholder
#    (<w:b/><w:bCs/>). First clone a fresh blank paragraph after it (so the
#    following placeholders are untouched), THEN overwrite the original
#    placeholder with the "16:30-18:30 ... página IA" text (bookmarked) and
#    the newly cloned one with "19:30-20:00 ... Código." text.
# ---------------------------------------------------------------------------
$idxPlaceholder = $idx832 + 1
$placeholderRange = $d.Paragraphs.Item($idxPlaceholder).Range
$placeholderRange.InsertParagraphAfter() | Out-Null

$body133 = @'
<w:p>
  <w:pPr>
    <w:ind w:left="360"/>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>16:30-</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">18:30 </w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_Hlk135764864"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Creación del contenido y diseño de la página IA</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>: imágenes y textos.</w:t>
  </w:r>
  <w:bookmarkEnd w:id="0"/>
</w:p>
'@

$range133 = $d.Paragraphs.Item($idxPlaceholder).Range
$range133.InsertXML($pkgHeader + $body133 + $pkgFooter) | Out-Null

$idxNew = $idxPlaceholder + 1

$body134 = @'
<w:p>
  <w:pPr>
    <w:ind w:left="360"/>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:u w:val="single"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">19:30-20:00 </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">Creación del contenido y diseño de la página IA: </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Código</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>.</w:t>
  </w:r>
</w:p>
'@

$range134 = $d.Paragraphs.Item($idxNew).Range
$range134.InsertXML($pkgHeader + $body134 + $pkgFooter) | Out-Null

Write-Output "Edit applied successfully."
